$d = $word.ActiveDocument

$replacements = @(
    @{old="719÷7=102, 5"; new="295÷7=42, 1"},
    @{old="616÷8=77, 0"; new="221÷4=55, 1"},
    @{old="919÷2=459, 1"; new="354÷4=88, 2"},
    @{old="820÷5=164, 0"; new="951÷2=475, 1"},
    @{old="232÷4=58, 0"; new="425÷2=212, 1"},
    @{old="395÷9=43, 8"; new="935÷9=103, 8"},
    @{old="251÷2=125, 1"; new="978÷9=108, 6"},
    @{old="304÷2=152, 0"; new="976÷4=244, 0"},
    @{old="708÷4=177, 0"; new="237÷7=33, 6"},
    @{old="185÷9=20, 5"; new="686÷8=85, 6"},
    @{old="678÷6=113, 0"; new="985÷8=123, 1"},
    @{old="196÷2=98, 0"; new="200÷6=33, 2"},
    @{old="457÷2=228, 1"; new="643÷2=321, 1"},
    @{old="654÷4=163, 2"; new="841÷7=120, 1"},
    @{old="187÷2=93, 1"; new="673÷2=336, 1"},
    @{old="269÷5=53, 4"; new="524÷6=87, 2"},
    @{old="656÷5=131, 1"; new="876÷5=175, 1"},
    @{old="468÷7=66, 6"; new="850÷4=212, 2"},
    @{old="792÷9=88, 0"; new="838÷9=93, 1"},
    @{old="308÷9=34, 2"; new="492÷4=123, 0"},
    @{old="867÷6=144, 3"; new="923÷6=153, 5"},
    @{old="348÷2=174, 0"; new="912÷8=114, 0"},
    @{old="200÷7=28, 4"; new="842÷2=421, 0"},
    @{old="574÷9=63, 7"; new="940÷6=156, 4"},
    @{old="610÷3=203, 1"; new="305÷3=101, 2"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Output "Done: $($replacements.Count) replacements applied"
